$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the sample number text "E7420" -> "E7420L" (shared across G2:G37)
$ws.Range("G2:G37").Value = "E7420L"

# 2. Replace the FALSE() formulas in H2:H37 with a literal boolean FALSE value
$ws.Range("H2:H37").Value = $false

# 3. Scroll the view so that A5 is the top-left visible cell
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
